$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The crawl timestamp column (O) was refreshed for every data row (rows 2-64)
# from "2022-09-04 07:02:22" to "2022-09-04 20:58:05".
for ($row = 2; $row -le 64; $row++) {
    $ws.Cells.Item($row, 15).Value = "2022-09-04 20:58:05"
}
